$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 170-171), pushing
# the existing rows 170-195 down to 172-197.
$ws.Rows("170:171").Insert()

# --- New row 170: weekly update, "Primera" ---
$ws.Range("A170").Value = 1
$ws.Range("B170").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C170").Value = "Arica y Parinacota"
$ws.Range("D170").Value = 44946
$ws.Range("E170").Value = 15
$ws.Range("F170").Value = "Fruta"
$ws.Range("G170").Value = 100108
$ws.Range("H170").Value = "Tropicales y subtropicales"
$ws.Range("I170").Value = 100108002
$ws.Range("J170").Value = "Mango"
$ws.Range("K170").Value = "Sin especificar"
$ws.Range("L170").Value = "Primera"
$ws.Range("M170").Value = 900
$ws.Range("N170").Value = 4500
$ws.Range("O170").Value = 5000
$ws.Range("P170").Value = 4714
$ws.Range("Q170").Value = "$/bandeja 4 kilos"
$ws.Range("R170").Value = "Perú"
$ws.Range("S170").Value = 1178
$ws.Range("T170").Value = 4

# --- New row 171: weekly update, "Segunda" ---
$ws.Range("A171").Value = 1
$ws.Range("B171").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C171").Value = "Arica y Parinacota"
$ws.Range("D171").Value = 44946
$ws.Range("E171").Value = 15
$ws.Range("F171").Value = "Fruta"
$ws.Range("G171").Value = 100108
$ws.Range("H171").Value = "Tropicales y subtropicales"
$ws.Range("I171").Value = 100108002
$ws.Range("J171").Value = "Mango"
$ws.Range("K171").Value = "Sin especificar"
$ws.Range("L171").Value = "Segunda"
$ws.Range("M171").Value = 700
$ws.Range("N171").Value = 4500
$ws.Range("O171").Value = 5000
$ws.Range("P171").Value = 4857
$ws.Range("Q171").Value = "$/bandeja 4 kilos"
$ws.Range("R171").Value = "Perú"
$ws.Range("S171").Value = 1214
$ws.Range("T171").Value = 4
